# Applies the "Updated symbol list" data refresh to the crypto price sheet.
# Column D holds numeric-looking values stored as TEXT (inlineStr in the
# original OOXML); writing a plain numeric-looking string makes Excel
# coerce it to a Number, so those assignments use a leading apostrophe to
# force text, then reset Style to "Normal" so no stray quote-prefix number
# format/style is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    $rng = $ws.Range($rangeAddress)
    $rng.Value = "'" + $value
    $rng.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2"  "241.20"
Set-TextValue "D3"  "22.34"
Set-TextValue "D4"  "5.521"
Set-TextValue "D5"  "0.05586"
Set-TextValue "D7"  "6.470"
Set-TextValue "D8"  "1.078"
Set-TextValue "D9"  "0.8033"
Set-TextValue "D10" "0.1419"
Set-TextValue "D11" "0.07398"
Set-TextValue "D12" "0.03256"
Set-TextValue "D13" "0.02996"
Set-TextValue "D14" "0.09250"
Set-TextValue "D15" "0.001674"
Set-TextValue "D16" "3.263"
Set-TextValue "D17" "0.04717"
Set-TextValue "D21" "0.003801"
Set-TextValue "D23" "0.0004777"
Set-TextValue "D25" "2.133"
Set-TextValue "D27" "0.1305"
Set-TextValue "D40" "0.04173"
Set-TextValue "D41" "0.006989"
Set-TextValue "D42" "0.003499"
Set-TextValue "D43" "0.1043"
Set-TextValue "D44" "0.009179"
Set-TextValue "D45" "0.00005491"
Set-TextValue "D48" "0.03038"

# --- Volume(1h) label (column E) updates ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E40").Value = "39IDEXIDEX"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("E48").Value = "47BOLOBOLO"

# --- Row 41 / 43 swap: BKEXToken <-> KickToken (Coin name + Link) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
